$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "64.034.02"
$ws.Range("E2").Value = "  -2.44%  "
# Row 3
Set-TextValue "D3" "3.343.60"
$ws.Range("E3").Value = "  -4.10%  "
# Row 4
$ws.Range("E4").Value = "  -0.17%  "
# Row 5
Set-TextValue "D5" "558.42"
$ws.Range("E5").Value = "  -0.89%  "
# Row 6
Set-TextValue "D6" "176.94"
$ws.Range("E6").Value = "  -0.17%  "
# Row 7
Set-TextValue "D7" "0.616"
$ws.Range("E7").Value = "  -2.01%  "
# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue "D8" "3.336.98"
$ws.Range("E8").Value = "  -4.27%  "
# Row 9
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D9" "1.00"
$ws.Range("E9").Value = "  -0.11%  "
# Row 10
Set-TextValue "D10" "0.623"
$ws.Range("E10").Value = "  -0.65%  "
# Row 11
Set-TextValue "D11" "0.160"
$ws.Range("E11").Value = "  +4.06%  "
# Row 12
Set-TextValue "D12" "53.71"
$ws.Range("E12").Value = "  -0.59%  "
# Row 13
Set-TextValue "D13" "0.0000269"
$ws.Range("E13").Value = "  -0.02%  "
# Row 14
Set-TextValue "D14" "9.06"
$ws.Range("E14").Value = "  -0.65%  "
# Row 15
Set-TextValue "D15" "3.872.56"
$ws.Range("E15").Value = "  -4.29%  "
# Row 16
Set-TextValue "D16" "18.37"
$ws.Range("E16").Value = "  +0.93%  "
# Row 17
Set-TextValue "D17" "0.119"
$ws.Range("E17").Value = "  -1.82%  "
# Row 18
Set-TextValue "D18" "3.343.07"
$ws.Range("E18").Value = "  -4.27%  "
# Row 19
Set-TextValue "D19" "11.85"
$ws.Range("E19").Value = "  -1.53%  "
# Row 20
Set-TextValue "D20" "63.934.25"
$ws.Range("E20").Value = "  -2.69%  "
# Row 21
Set-TextValue "D21" "0.983"
$ws.Range("E21").Value = "  -0.98%  "
# Row 22
Set-TextValue "D22" "447.77"
$ws.Range("E22").Value = "  +8.95%  "
# Row 23
Set-TextValue "D23" "4.55"
$ws.Range("E23").Value = "  +10.89%  "
# Row 24
Set-TextValue "D24" "4.11"
$ws.Range("E24").Value = "  +0.23%  "
# Row 25
Set-TextValue "D25" "84.91"
$ws.Range("E25").Value = "  +0.07%  "
# Row 26
Set-TextValue "D26" "13.24"
$ws.Range("E26").Value = "  +4.02%  "
# Row 27
Set-TextValue "D27" "10.69"
$ws.Range("E27").Value = "  -1.41%  "
# Row 28
Set-TextValue "D28" "2.83"
$ws.Range("E28").Value = "  +0.19%  "
# Row 29
Set-TextValue "D29" "8.77"
$ws.Range("E29").Value = "  -1.52%  "
# Row 30
Set-TextValue "D30" "29.67"
$ws.Range("E30").Value = "  -1.11%  "
# Row 31
Set-TextValue "D31" "6.59"
$ws.Range("E31").Value = "  +4.39%  "
# Row 32
Set-TextValue "D32" "588.74"
$ws.Range("E32").Value = "  -4.86%  "
# Row 33
Set-TextValue "D33" "11.49"
$ws.Range("E33").Value = "  -0.75%  "
# Row 34
Set-TextValue "D34" "0.107"
$ws.Range("E34").Value = "  -1.08%  "
# Row 35
Set-TextValue "D35" "58.69"
$ws.Range("E35").Value = "  -0.44%  "
# Row 36
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.10%  "
# Row 37
$ws.Range("E37").Value = "  -3.14%  "
# Row 38
Set-TextValue "D38" "3.54"
$ws.Range("E38").Value = "  +4.36%  "
# Row 39
Set-TextValue "D39" "35.60"
$ws.Range("E39").Value = "  -3.30%  "
# Row 40
Set-TextValue "D40" "0.0₃0750"
$ws.Range("E40").Value = "  -4.63%  "
# Row 41
Set-TextValue "D41" "0.368"
$ws.Range("E41").Value = "  -2.00%  "
# Row 42
Set-TextValue "D42" "3.128.77"
$ws.Range("E42").Value = "  -6.95%  "
# Row 43
$ws.Range("E43").Value = "  -0.43%  "
# Row 44
Set-TextValue "D44" "2.87"
$ws.Range("E44").Value = "  +0.06%  "
# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0406"
$ws.Range("E45").Value = "  -1.74%  "
# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D46" "3.19"
$ws.Range("E46").Value = "  -1.92%  "
# Row 47
Set-TextValue "D47" "2.45"
$ws.Range("E47").Value = "  -1.50%  "
# Row 48
$ws.Range("E48").Value = "  -1.53%  "
# Row 49
$ws.Range("E49").Value = "  -4.30%  "
# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D50" "8.19"
$ws.Range("E50").Value = "  -2.16%  "
# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "133.51"
$ws.Range("E51").Value = "  -2.73%  "
